# "major cleaning of unnessecary stuff"
#
# This workbook is a chain of test-registration fixture sheets:
#   Sheet4 (raw data entry) --> Sheet2/Sheet5 (formulas referencing Sheet4)
#   --> Sheet3/Sheet6 (formulas referencing Sheet2/Sheet5)
#   --> Sheet7/Sheet8 (formulas referencing Sheet5/Sheet6)
# Editing the handful of "source of truth" cells below lets every
# downstream formula cell recalc on its own.

$wb = $excel.ActiveWorkbook

# --- Reg sheet: bump the two standalone test-account labels ---
$wsReg = $wb.Worksheets.Item("Reg")
$wsReg.Range("B3").Value = "tavalinetont24"
$wsReg.Range("B4").Value = "puhtaloom24"

# --- Sheet4: the master row every other sheet's formulas chase ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A2").Value = "tavalinetont42"
$ws4.Range("C2").Value = "puhtaloom42"
$ws4.Range("E2").Value = "filmweird42"

# --- Sheet5: bump the mailinator addresses to the next run ---
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("I2").Value = "selentest23@mailinator.com"
$ws5.Range("I3").Value = "testimeauto23@mailinator.com"
$ws5.Range("I4").Value = "vahekonto23@mailinator.com"

# --- Sheet8: small wording fixes (drop the accented characters) ---
$ws8 = $wb.Worksheets.Item("Sheet8")
$ws8.Range("F2").Value = "Eks me koik oleme natuke imelikud"
$ws8.Range("E4").Value = "Teda ei hairi miski"

# --- Sheet9: fill in the two previously-missing upload sizes ---
$ws9 = $wb.Worksheets.Item("Sheet9")
$ws9.Range("G3").NumberFormat = "@"
$ws9.Range("G3").Value = "1523"
$ws9.Range("H3").NumberFormat = "@"
$ws9.Range("H3").Value = "1202"

# --- restore each sheet's recorded selection / active cell ---
$wsReg.Activate()
$wsReg.Range("B3").Select()

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("E11").Select()

$ws4.Activate()
$ws4.Range("E2").Select()

$ws5.Activate()
$ws5.Range("A4").Select()

$ws8.Activate()
$ws8.Range("F4").Select()

# Sheet4 is the tab that was active/selected in the saved workbook.
$ws4.Activate()
